$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update realized (N) / remaining (O) contract value figures ----
# Row 7
$ws.Cells.Item(7, 14).Value2 = 522500000
$ws.Cells.Item(7, 15).Value2 = 47500000
# Row 9
$ws.Cells.Item(9, 14).Value2 = 2090992550
$ws.Cells.Item(9, 15).Value2 = 297915680
# Row 43
$ws.Cells.Item(43, 14).Value2 = 1297948045
$ws.Cells.Item(43, 15).Value2 = 262639128
# Row 49
$ws.Cells.Item(49, 14).Value2 = 395250000
$ws.Cells.Item(49, 15).Value2 = 0
# Row 50
$ws.Cells.Item(50, 14).Value2 = 118010250
$ws.Cells.Item(50, 15).Value2 = 37153750
# Row 65
$ws.Cells.Item(65, 14).Value2 = 98346000
$ws.Cells.Item(65, 15).Value2 = 0
# Row 66
$ws.Cells.Item(66, 14).Value2 = 86802000
$ws.Cells.Item(66, 15).Value2 = 0
# Row 67
$ws.Cells.Item(67, 14).Value2 = 99900000
$ws.Cells.Item(67, 15).Value2 = 0
# Row 68
$ws.Cells.Item(68, 14).Value2 = 191440000
$ws.Cells.Item(68, 15).Value2 = 0
# Row 69
$ws.Cells.Item(69, 14).Value2 = 146381000
$ws.Cells.Item(69, 15).Value2 = 0

# ---- New contract rows (70-85): SP2D LS Non Kontraktual entries ----
$ws.Cells.Item(70, 1).Value2 = 66
$ws.Cells.Item(70, 2).Value2 = 626402
$ws.Cells.Item(70, 3).NumberFormat = "@"
$ws.Cells.Item(70, 3).Value2 = 'A/175.22009480/0/0'
$ws.Cells.Item(70, 4).NumberFormat = "@"
$ws.Cells.Item(70, 4).Value2 = 'PT. ASTRIDO JAYA MOBILINDO'
$ws.Cells.Item(70, 5).NumberFormat = "@"
$ws.Cells.Item(70, 5).Value2 = '25-NOV-22'
$ws.Cells.Item(70, 6).NumberFormat = "@"
$ws.Cells.Item(70, 6).Value2 = '25-NOV-22'
$ws.Cells.Item(70, 7).NumberFormat = "@"
$ws.Cells.Item(70, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(70, 8).NumberFormat = "@"
$ws.Cells.Item(70, 8).Value2 = 'SPK.4651/PPK/BRSDM.01/XI/2022'
$ws.Cells.Item(70, 9).NumberFormat = "@"
$ws.Cells.Item(70, 9).Value2 = 'Pengadaan Kendaraan Dinas Pimpinan pada Sekretariat BRSDM'
$ws.Cells.Item(70, 10).Value2 = 532111
$ws.Cells.Item(70, 11).NumberFormat = "@"
$ws.Cells.Item(70, 11).Value2 = '25-NOV-22'
$ws.Cells.Item(70, 12).NumberFormat = "@"
$ws.Cells.Item(70, 12).Value2 = '''''01-DEC-22'
$ws.Cells.Item(70, 13).Value2 = 712785000
$ws.Cells.Item(70, 14).Value2 = 712785000
$ws.Cells.Item(70, 15).Value2 = 0

$ws.Cells.Item(71, 1).Value2 = 67
$ws.Cells.Item(71, 2).Value2 = 626402
$ws.Cells.Item(71, 3).NumberFormat = "@"
$ws.Cells.Item(71, 3).Value2 = 'A/175.22009496/0/0'
$ws.Cells.Item(71, 4).NumberFormat = "@"
$ws.Cells.Item(71, 4).Value2 = 'UPP-TPA KOPERASI PERENCANAAN'
$ws.Cells.Item(71, 5).NumberFormat = "@"
$ws.Cells.Item(71, 5).Value2 = '25-NOV-22'
$ws.Cells.Item(71, 6).NumberFormat = "@"
$ws.Cells.Item(71, 6).Value2 = '18-NOV-22'
$ws.Cells.Item(71, 7).NumberFormat = "@"
$ws.Cells.Item(71, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(71, 8).NumberFormat = "@"
$ws.Cells.Item(71, 8).Value2 = '480/PPK.PUSDIK/PL.430/XI/2022'
$ws.Cells.Item(71, 9).NumberFormat = "@"
$ws.Cells.Item(71, 9).Value2 = 'Pelaksanaan Tes TOEFL dan TPA'
$ws.Cells.Item(71, 10).Value2 = 522191
$ws.Cells.Item(71, 11).NumberFormat = "@"
$ws.Cells.Item(71, 11).Value2 = '21-NOV-22'
$ws.Cells.Item(71, 12).NumberFormat = "@"
$ws.Cells.Item(71, 12).Value2 = '22-NOV-22'
$ws.Cells.Item(71, 13).Value2 = 139805500
$ws.Cells.Item(71, 14).Value2 = 139805500
$ws.Cells.Item(71, 15).Value2 = 0

$ws.Cells.Item(72, 1).Value2 = 68
$ws.Cells.Item(72, 2).Value2 = 626402
$ws.Cells.Item(72, 3).NumberFormat = "@"
$ws.Cells.Item(72, 3).Value2 = 'A/175.22009544/0/0'
$ws.Cells.Item(72, 4).NumberFormat = "@"
$ws.Cells.Item(72, 4).Value2 = 'PT. SAKHA PRATAMA MANDIRI'
$ws.Cells.Item(72, 5).NumberFormat = "@"
$ws.Cells.Item(72, 5).Value2 = '28-NOV-22'
$ws.Cells.Item(72, 6).NumberFormat = "@"
$ws.Cells.Item(72, 6).Value2 = '24-NOV-22'
$ws.Cells.Item(72, 7).NumberFormat = "@"
$ws.Cells.Item(72, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(72, 8).NumberFormat = "@"
$ws.Cells.Item(72, 8).Value2 = 'SPK-4631/PPK/BRSDM.1/XI/2022'
$ws.Cells.Item(72, 9).NumberFormat = "@"
$ws.Cells.Item(72, 9).Value2 = 'Pekerjaan Jasa Perbaikan/Perawatan Lift dan Komponennya pada Sekretariat BRSDM'
$ws.Cells.Item(72, 10).Value2 = 523121
$ws.Cells.Item(72, 11).NumberFormat = "@"
$ws.Cells.Item(72, 11).Value2 = '24-NOV-22'
$ws.Cells.Item(72, 12).NumberFormat = "@"
$ws.Cells.Item(72, 12).Value2 = '15-DEC-22'
$ws.Cells.Item(72, 13).Value2 = 98000000
$ws.Cells.Item(72, 14).Value2 = 0
$ws.Cells.Item(72, 15).Value2 = 98000000

$ws.Cells.Item(73, 1).Value2 = 69
$ws.Cells.Item(73, 2).Value2 = 626402
$ws.Cells.Item(73, 3).NumberFormat = "@"
$ws.Cells.Item(73, 3).Value2 = 'A/175.22009792/0/0'
$ws.Cells.Item(73, 4).NumberFormat = "@"
$ws.Cells.Item(73, 4).Value2 = 'PT ALFABET INDO KREATIF'
$ws.Cells.Item(73, 5).NumberFormat = "@"
$ws.Cells.Item(73, 5).Value2 = '''''05-DEC-22'
$ws.Cells.Item(73, 6).NumberFormat = "@"
$ws.Cells.Item(73, 6).Value2 = '30-NOV-22'
$ws.Cells.Item(73, 7).NumberFormat = "@"
$ws.Cells.Item(73, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(73, 8).NumberFormat = "@"
$ws.Cells.Item(73, 8).Value2 = '507/PPK.PUSDIK/PL.430/XI/2022'
$ws.Cells.Item(73, 9).NumberFormat = "@"
$ws.Cells.Item(73, 9).Value2 = 'Pengadaan Kalender Pusat Pendidikan KP'
$ws.Cells.Item(73, 10).Value2 = 521111
$ws.Cells.Item(73, 11).NumberFormat = "@"
$ws.Cells.Item(73, 11).Value2 = '30-NOV-22'
$ws.Cells.Item(73, 12).NumberFormat = "@"
$ws.Cells.Item(73, 12).Value2 = '14-DEC-22'
$ws.Cells.Item(73, 13).Value2 = 62770500
$ws.Cells.Item(73, 14).Value2 = 0
$ws.Cells.Item(73, 15).Value2 = 62770500

$ws.Cells.Item(74, 1).Value2 = 70
$ws.Cells.Item(74, 2).Value2 = 626402
$ws.Cells.Item(74, 3).NumberFormat = "@"
$ws.Cells.Item(74, 3).Value2 = 'A/175.22009873/0/0'
$ws.Cells.Item(74, 4).NumberFormat = "@"
$ws.Cells.Item(74, 4).Value2 = 'DIN LAW GROUP'
$ws.Cells.Item(74, 5).NumberFormat = "@"
$ws.Cells.Item(74, 5).Value2 = '''''07-DEC-22'
$ws.Cells.Item(74, 6).NumberFormat = "@"
$ws.Cells.Item(74, 6).Value2 = '''''01-DEC-22'
$ws.Cells.Item(74, 7).NumberFormat = "@"
$ws.Cells.Item(74, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(74, 8).NumberFormat = "@"
$ws.Cells.Item(74, 8).Value2 = 'SPK-4533/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(74, 9).NumberFormat = "@"
$ws.Cells.Item(74, 9).Value2 = 'Pekerjaan Jasa Konsultasi Hukum Kegiatan Prioritas BRSDM'
$ws.Cells.Item(74, 10).Value2 = 522191
$ws.Cells.Item(74, 11).NumberFormat = "@"
$ws.Cells.Item(74, 11).Value2 = '''''01-DEC-22'
$ws.Cells.Item(74, 12).NumberFormat = "@"
$ws.Cells.Item(74, 12).Value2 = '19-DEC-22'
$ws.Cells.Item(74, 13).Value2 = 176000000
$ws.Cells.Item(74, 14).Value2 = 0
$ws.Cells.Item(74, 15).Value2 = 176000000

$ws.Cells.Item(75, 1).Value2 = 71
$ws.Cells.Item(75, 2).Value2 = 626402
$ws.Cells.Item(75, 3).NumberFormat = "@"
$ws.Cells.Item(75, 3).Value2 = 'A/175.22009874/0/0'
$ws.Cells.Item(75, 4).NumberFormat = "@"
$ws.Cells.Item(75, 4).Value2 = 'PT. PROSPERITAS FORTUNA INDONESIA'
$ws.Cells.Item(75, 5).NumberFormat = "@"
$ws.Cells.Item(75, 5).Value2 = '''''07-DEC-22'
$ws.Cells.Item(75, 6).NumberFormat = "@"
$ws.Cells.Item(75, 6).Value2 = '''''02-DEC-22'
$ws.Cells.Item(75, 7).NumberFormat = "@"
$ws.Cells.Item(75, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(75, 8).NumberFormat = "@"
$ws.Cells.Item(75, 8).Value2 = 'SPK-4652/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(75, 9).NumberFormat = "@"
$ws.Cells.Item(75, 9).Value2 = 'Pekerjaan Jasa Rehab Ruang Rapat Sekretariat BRSDM'
$ws.Cells.Item(75, 10).Value2 = 523111
$ws.Cells.Item(75, 11).NumberFormat = "@"
$ws.Cells.Item(75, 11).Value2 = '''''02-DEC-22'
$ws.Cells.Item(75, 12).NumberFormat = "@"
$ws.Cells.Item(75, 12).Value2 = '12-DEC-22'
$ws.Cells.Item(75, 13).Value2 = 89982000
$ws.Cells.Item(75, 14).Value2 = 0
$ws.Cells.Item(75, 15).Value2 = 89982000

$ws.Cells.Item(76, 1).Value2 = 72
$ws.Cells.Item(76, 2).Value2 = 626402
$ws.Cells.Item(76, 3).NumberFormat = "@"
$ws.Cells.Item(76, 3).Value2 = 'A/175.22009996/0/0'
$ws.Cells.Item(76, 4).NumberFormat = "@"
$ws.Cells.Item(76, 4).Value2 = 'CV. TECHNO ENERGY'
$ws.Cells.Item(76, 5).NumberFormat = "@"
$ws.Cells.Item(76, 5).Value2 = '12-DEC-22'
$ws.Cells.Item(76, 6).NumberFormat = "@"
$ws.Cells.Item(76, 6).Value2 = '''''08-DEC-22'
$ws.Cells.Item(76, 7).NumberFormat = "@"
$ws.Cells.Item(76, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(76, 8).NumberFormat = "@"
$ws.Cells.Item(76, 8).Value2 = 'SPK_4831/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(76, 9).NumberFormat = "@"
$ws.Cells.Item(76, 9).Value2 = 'Pekerjaan Pembuatan Kalender BRSDM 2023'
$ws.Cells.Item(76, 10).Value2 = 521111
$ws.Cells.Item(76, 11).NumberFormat = "@"
$ws.Cells.Item(76, 11).Value2 = '''''08-DEC-22'
$ws.Cells.Item(76, 12).NumberFormat = "@"
$ws.Cells.Item(76, 12).Value2 = '15-DEC-22'
$ws.Cells.Item(76, 13).Value2 = 109057000
$ws.Cells.Item(76, 14).Value2 = 0
$ws.Cells.Item(76, 15).Value2 = 109057000

$ws.Cells.Item(77, 1).Value2 = 73
$ws.Cells.Item(77, 2).Value2 = 626402
$ws.Cells.Item(77, 3).NumberFormat = "@"
$ws.Cells.Item(77, 3).Value2 = 'A/175.22010016/0/0'
$ws.Cells.Item(77, 4).NumberFormat = "@"
$ws.Cells.Item(77, 4).Value2 = 'MIRAH SEGAR'
$ws.Cells.Item(77, 5).NumberFormat = "@"
$ws.Cells.Item(77, 5).Value2 = '12-DEC-22'
$ws.Cells.Item(77, 6).NumberFormat = "@"
$ws.Cells.Item(77, 6).Value2 = '''''07-DEC-22'
$ws.Cells.Item(77, 7).NumberFormat = "@"
$ws.Cells.Item(77, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(77, 8).NumberFormat = "@"
$ws.Cells.Item(77, 8).Value2 = '522/PPK.PUSDIK/PL.430/XII/2022'
$ws.Cells.Item(77, 9).NumberFormat = "@"
$ws.Cells.Item(77, 9).Value2 = 'Pelaksanaan Fullboard Meeting Kegiatan Rapat Koordinasi Tugas Belajar dan Izin Belajar Tahun 2022'
$ws.Cells.Item(77, 10).Value2 = 524119
$ws.Cells.Item(77, 11).NumberFormat = "@"
$ws.Cells.Item(77, 11).Value2 = '''''07-DEC-22'
$ws.Cells.Item(77, 12).NumberFormat = "@"
$ws.Cells.Item(77, 12).Value2 = '''''09-DEC-22'
$ws.Cells.Item(77, 13).Value2 = 110250000
$ws.Cells.Item(77, 14).Value2 = 0
$ws.Cells.Item(77, 15).Value2 = 110250000

$ws.Cells.Item(78, 1).Value2 = 74
$ws.Cells.Item(78, 2).Value2 = 626402
$ws.Cells.Item(78, 3).NumberFormat = "@"
$ws.Cells.Item(78, 3).Value2 = 'A/175.22010060/0/0'
$ws.Cells.Item(78, 4).NumberFormat = "@"
$ws.Cells.Item(78, 4).Value2 = 'PT. NOVAL INDO PRATAMA'
$ws.Cells.Item(78, 5).NumberFormat = "@"
$ws.Cells.Item(78, 5).Value2 = '13-DEC-22'
$ws.Cells.Item(78, 6).NumberFormat = "@"
$ws.Cells.Item(78, 6).Value2 = '''''08-DEC-22'
$ws.Cells.Item(78, 7).NumberFormat = "@"
$ws.Cells.Item(78, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(78, 8).NumberFormat = "@"
$ws.Cells.Item(78, 8).Value2 = 'SPK-4832/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(78, 9).NumberFormat = "@"
$ws.Cells.Item(78, 9).Value2 = 'Jasa Rehab Ruang Selasar Kepala BRSDM'
$ws.Cells.Item(78, 10).Value2 = 523111
$ws.Cells.Item(78, 11).NumberFormat = "@"
$ws.Cells.Item(78, 11).Value2 = '''''08-DEC-22'
$ws.Cells.Item(78, 12).NumberFormat = "@"
$ws.Cells.Item(78, 12).Value2 = '12-DEC-22'
$ws.Cells.Item(78, 13).Value2 = 193150000
$ws.Cells.Item(78, 14).Value2 = 0
$ws.Cells.Item(78, 15).Value2 = 193150000

$ws.Cells.Item(79, 1).Value2 = 75
$ws.Cells.Item(79, 2).Value2 = 626402
$ws.Cells.Item(79, 3).NumberFormat = "@"
$ws.Cells.Item(79, 3).Value2 = 'A/175.22010136/0/0'
$ws.Cells.Item(79, 4).NumberFormat = "@"
$ws.Cells.Item(79, 4).Value2 = 'PT. ARGA SINAR TERANG'
$ws.Cells.Item(79, 5).NumberFormat = "@"
$ws.Cells.Item(79, 5).Value2 = '15-DEC-22'
$ws.Cells.Item(79, 6).NumberFormat = "@"
$ws.Cells.Item(79, 6).Value2 = '''''09-DEC-22'
$ws.Cells.Item(79, 7).NumberFormat = "@"
$ws.Cells.Item(79, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(79, 8).NumberFormat = "@"
$ws.Cells.Item(79, 8).Value2 = 'SPK-4852/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(79, 9).NumberFormat = "@"
$ws.Cells.Item(79, 9).Value2 = 'Jasa Instalasi Listrik dan AC Ruang Rapat BRSDM'
$ws.Cells.Item(79, 10).Value2 = 523121
$ws.Cells.Item(79, 11).NumberFormat = "@"
$ws.Cells.Item(79, 11).Value2 = '''''09-DEC-22'
$ws.Cells.Item(79, 12).NumberFormat = "@"
$ws.Cells.Item(79, 12).Value2 = '14-DEC-22'
$ws.Cells.Item(79, 13).Value2 = 196656000
$ws.Cells.Item(79, 14).Value2 = 0
$ws.Cells.Item(79, 15).Value2 = 196656000

$ws.Cells.Item(80, 1).Value2 = 76
$ws.Cells.Item(80, 2).Value2 = 626402
$ws.Cells.Item(80, 3).NumberFormat = "@"
$ws.Cells.Item(80, 3).Value2 = 'A/175.22010137/0/0'
$ws.Cells.Item(80, 4).NumberFormat = "@"
$ws.Cells.Item(80, 4).Value2 = 'PT. ARTA ANUGRAH SEJAHTERA'
$ws.Cells.Item(80, 5).NumberFormat = "@"
$ws.Cells.Item(80, 5).Value2 = '15-DEC-22'
$ws.Cells.Item(80, 6).NumberFormat = "@"
$ws.Cells.Item(80, 6).Value2 = '''''09-DEC-22'
$ws.Cells.Item(80, 7).NumberFormat = "@"
$ws.Cells.Item(80, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(80, 8).NumberFormat = "@"
$ws.Cells.Item(80, 8).Value2 = 'SPK-4853/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(80, 9).NumberFormat = "@"
$ws.Cells.Item(80, 9).Value2 = 'Jasa Pembuatan Dinding Partisi Area Ruang Rapat BRSDM'
$ws.Cells.Item(80, 10).Value2 = 523111
$ws.Cells.Item(80, 11).NumberFormat = "@"
$ws.Cells.Item(80, 11).Value2 = '''''09-DEC-22'
$ws.Cells.Item(80, 12).NumberFormat = "@"
$ws.Cells.Item(80, 12).Value2 = '14-DEC-22'
$ws.Cells.Item(80, 13).Value2 = 198918000
$ws.Cells.Item(80, 14).Value2 = 0
$ws.Cells.Item(80, 15).Value2 = 198918000

$ws.Cells.Item(81, 1).Value2 = 77
$ws.Cells.Item(81, 2).Value2 = 626402
$ws.Cells.Item(81, 3).NumberFormat = "@"
$ws.Cells.Item(81, 3).Value2 = 'A/175.22010138/0/0'
$ws.Cells.Item(81, 4).NumberFormat = "@"
$ws.Cells.Item(81, 4).Value2 = 'PT. AULIA BERLIAN KONSTRUKSI'
$ws.Cells.Item(81, 5).NumberFormat = "@"
$ws.Cells.Item(81, 5).Value2 = '15-DEC-22'
$ws.Cells.Item(81, 6).NumberFormat = "@"
$ws.Cells.Item(81, 6).Value2 = '''''09-DEC-22'
$ws.Cells.Item(81, 7).NumberFormat = "@"
$ws.Cells.Item(81, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(81, 8).NumberFormat = "@"
$ws.Cells.Item(81, 8).Value2 = 'SPK-4854/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(81, 9).NumberFormat = "@"
$ws.Cells.Item(81, 9).Value2 = 'Jasa Rehab Lantai dan Plafon Area Ruang Rapat BRSDM'
$ws.Cells.Item(81, 10).Value2 = 523111
$ws.Cells.Item(81, 11).NumberFormat = "@"
$ws.Cells.Item(81, 11).Value2 = '''''09-DEC-22'
$ws.Cells.Item(81, 12).NumberFormat = "@"
$ws.Cells.Item(81, 12).Value2 = '14-DEC-22'
$ws.Cells.Item(81, 13).Value2 = 162456659
$ws.Cells.Item(81, 14).Value2 = 0
$ws.Cells.Item(81, 15).Value2 = 162456659

$ws.Cells.Item(82, 1).Value2 = 78
$ws.Cells.Item(82, 2).Value2 = 626402
$ws.Cells.Item(82, 3).NumberFormat = "@"
$ws.Cells.Item(82, 3).Value2 = 'A/175.22010177/0/0'
$ws.Cells.Item(82, 4).NumberFormat = "@"
$ws.Cells.Item(82, 4).Value2 = 'PT. AULIA BERLIAN KONSTRUKSI'
$ws.Cells.Item(82, 5).NumberFormat = "@"
$ws.Cells.Item(82, 5).Value2 = '16-DEC-22'
$ws.Cells.Item(82, 6).NumberFormat = "@"
$ws.Cells.Item(82, 6).Value2 = '15-DEC-22'
$ws.Cells.Item(82, 7).NumberFormat = "@"
$ws.Cells.Item(82, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(82, 8).NumberFormat = "@"
$ws.Cells.Item(82, 8).Value2 = 'SPK-4936/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(82, 9).NumberFormat = "@"
$ws.Cells.Item(82, 9).Value2 = 'Jasa Pembuatan WLL Panel Ruang Rapat BRSDM'
$ws.Cells.Item(82, 10).Value2 = 523111
$ws.Cells.Item(82, 11).NumberFormat = "@"
$ws.Cells.Item(82, 11).Value2 = '15-DEC-22'
$ws.Cells.Item(82, 12).NumberFormat = "@"
$ws.Cells.Item(82, 12).Value2 = '19-DEC-22'
$ws.Cells.Item(82, 13).Value2 = 191086500
$ws.Cells.Item(82, 14).Value2 = 0
$ws.Cells.Item(82, 15).Value2 = 191086500

$ws.Cells.Item(83, 1).Value2 = 79
$ws.Cells.Item(83, 2).Value2 = 626402
$ws.Cells.Item(83, 3).NumberFormat = "@"
$ws.Cells.Item(83, 3).Value2 = 'A/175.22010178/0/0'
$ws.Cells.Item(83, 4).NumberFormat = "@"
$ws.Cells.Item(83, 4).Value2 = 'PT. ARTA ANUGRAH SEJAHTERA'
$ws.Cells.Item(83, 5).NumberFormat = "@"
$ws.Cells.Item(83, 5).Value2 = '16-DEC-22'
$ws.Cells.Item(83, 6).NumberFormat = "@"
$ws.Cells.Item(83, 6).Value2 = '15-DEC-22'
$ws.Cells.Item(83, 7).NumberFormat = "@"
$ws.Cells.Item(83, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(83, 8).NumberFormat = "@"
$ws.Cells.Item(83, 8).Value2 = 'SPK-4937/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(83, 9).NumberFormat = "@"
$ws.Cells.Item(83, 9).Value2 = 'Jasa Pembuatan Logo dan Custom Cabinet Ruang Rapat BRSDM'
$ws.Cells.Item(83, 10).Value2 = 523111
$ws.Cells.Item(83, 11).NumberFormat = "@"
$ws.Cells.Item(83, 11).Value2 = '15-DEC-22'
$ws.Cells.Item(83, 12).NumberFormat = "@"
$ws.Cells.Item(83, 12).Value2 = '19-DEC-22'
$ws.Cells.Item(83, 13).Value2 = 141192000
$ws.Cells.Item(83, 14).Value2 = 0
$ws.Cells.Item(83, 15).Value2 = 141192000

$ws.Cells.Item(84, 1).Value2 = 80
$ws.Cells.Item(84, 2).Value2 = 626402
$ws.Cells.Item(84, 3).NumberFormat = "@"
$ws.Cells.Item(84, 3).Value2 = 'A/175.22010179/0/0'
$ws.Cells.Item(84, 4).NumberFormat = "@"
$ws.Cells.Item(84, 4).Value2 = 'PT. ARGA SINAR TERANG'
$ws.Cells.Item(84, 5).NumberFormat = "@"
$ws.Cells.Item(84, 5).Value2 = '16-DEC-22'
$ws.Cells.Item(84, 6).NumberFormat = "@"
$ws.Cells.Item(84, 6).Value2 = '15-DEC-22'
$ws.Cells.Item(84, 7).NumberFormat = "@"
$ws.Cells.Item(84, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(84, 8).NumberFormat = "@"
$ws.Cells.Item(84, 8).Value2 = 'SPK-4935/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(84, 9).NumberFormat = "@"
$ws.Cells.Item(84, 9).Value2 = 'Pekerjaan Pintu dan Partisi Kaca Ruang Rapat BRSDM'
$ws.Cells.Item(84, 10).Value2 = 523111
$ws.Cells.Item(84, 11).NumberFormat = "@"
$ws.Cells.Item(84, 11).Value2 = '15-DEC-22'
$ws.Cells.Item(84, 12).NumberFormat = "@"
$ws.Cells.Item(84, 12).Value2 = '19-DEC-22'
$ws.Cells.Item(84, 13).Value2 = 130101990
$ws.Cells.Item(84, 14).Value2 = 0
$ws.Cells.Item(84, 15).Value2 = 130101990

$ws.Cells.Item(85, 1).Value2 = 81
$ws.Cells.Item(85, 2).Value2 = 626402
$ws.Cells.Item(85, 3).NumberFormat = "@"
$ws.Cells.Item(85, 3).Value2 = 'A/175.22010180/0/0'
$ws.Cells.Item(85, 4).NumberFormat = "@"
$ws.Cells.Item(85, 4).Value2 = 'PT. NOVAL INDO PRATAMA'
$ws.Cells.Item(85, 5).NumberFormat = "@"
$ws.Cells.Item(85, 5).Value2 = '16-DEC-22'
$ws.Cells.Item(85, 6).NumberFormat = "@"
$ws.Cells.Item(85, 6).Value2 = '13-DEC-22'
$ws.Cells.Item(85, 7).NumberFormat = "@"
$ws.Cells.Item(85, 7).Value2 = 'Tidak terlambat'
$ws.Cells.Item(85, 8).NumberFormat = "@"
$ws.Cells.Item(85, 8).Value2 = 'SPK-4893/PPK/BRSDM.1/XII/2022'
$ws.Cells.Item(85, 9).NumberFormat = "@"
$ws.Cells.Item(85, 9).Value2 = 'Jasa Konstruksi Audience Bench Ruang Rapat Utama'
$ws.Cells.Item(85, 10).Value2 = 523111
$ws.Cells.Item(85, 11).NumberFormat = "@"
$ws.Cells.Item(85, 11).Value2 = '13-DEC-22'
$ws.Cells.Item(85, 12).NumberFormat = "@"
$ws.Cells.Item(85, 12).Value2 = '19-DEC-22'
$ws.Cells.Item(85, 13).Value2 = 193584000
$ws.Cells.Item(85, 14).Value2 = 0
$ws.Cells.Item(85, 15).Value2 = 193584000

# ---- Clean up auto-applied "Text" / "quote prefix" cell styling ----
# Text columns above had NumberFormat forced to "@" (Text) before the
# assignment so that date-shaped values (e.g. "25-NOV-22") are kept as
# literal text instead of being auto-converted to Excel date serial
# numbers, and so that values starting with a literal apostrophe (which
# were doubled before assignment so the apostrophe survives as real text
# content instead of being consumed as Excel's "quote prefix" marker)
# keep that apostrophe. Afterwards, clear the cell formatting on the
# whole new range back to the default (unstyled) look used throughout
# the rest of the table.
$ws.Range("A70:O85").ClearFormats()
